$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Support Beam" renamed to "Wooden Plank" (row 13, Assets column)
$ws.Range("B13").Value = "Wooden Plank"

# Mark Building 1, Building 2 and Crane features as Done in the Status column
$ws.Range("G7").Value = "Done"
$ws.Range("G8").Value = "Done"
$ws.Range("G16").Value = "Done"
$ws.Range("G16").HorizontalAlignment = -4152

# Update current selection
$ws.Range("I13").Select()
